# Update the header date.
$d = $word.ActiveDocument
[void]$d.Content.Find.Execute("2025-07-29 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-07-30 Wednesday", 2)

# Update the answer table values in place (table shape / row count is unchanged —
# only the text content of the five "problem rows" shifts).
$t = $d.Tables.Item(1)

# Row 1
$t.Cell(1, 1).Range.Text = "41÷7=5, 6"
$t.Cell(1, 2).Range.Text = "49÷6=8, 1"
$t.Cell(1, 3).Range.Text = "99÷2=49, 1"
$t.Cell(1, 4).Range.Text = "53÷3=17, 2"
$t.Cell(1, 5).Range.Text = "17÷7=2, 3"

# Row 5
$t.Cell(5, 1).Range.Text = "64÷7=9, 1"
$t.Cell(5, 2).Range.Text = "90÷5=18, 0"
$t.Cell(5, 3).Range.Text = "58÷9=6, 4"
$t.Cell(5, 4).Range.Text = "71÷7=10, 1"
$t.Cell(5, 5).Range.Text = "15÷6=2, 3"

# Row 9
$t.Cell(9, 1).Range.Text = "65÷3=21, 2"
$t.Cell(9, 2).Range.Text = "34÷9=3, 7"
$t.Cell(9, 3).Range.Text = "92÷5=18, 2"
$t.Cell(9, 4).Range.Text = "40÷4=10, 0"
$t.Cell(9, 5).Range.Text = "98÷5=19, 3"

# Row 13
$t.Cell(13, 1).Range.Text = "27÷2=13, 1"
$t.Cell(13, 2).Range.Text = "25÷5=5, 0"
$t.Cell(13, 3).Range.Text = "50÷8=6, 2"
$t.Cell(13, 4).Range.Text = "22÷4=5, 2"
$t.Cell(13, 5).Range.Text = "18÷4=4, 2"

# Row 17
$t.Cell(17, 1).Range.Text = "87÷8=10, 7"
$t.Cell(17, 2).Range.Text = "70÷5=14, 0"
$t.Cell(17, 3).Range.Text = "26÷6=4, 2"
$t.Cell(17, 4).Range.Text = "24÷8=3, 0"
$t.Cell(17, 5).Range.Text = "79÷6=13, 1"
